$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 112388117
$ws.Range("B4").Value = 77650
$ws.Range("Q4").Value = 557810
$ws.Range("R4").Value = 7069645

# Row 5
$ws.Range("A5").Value = 112388107
$ws.Range("B5").Value = 81385
$ws.Range("E5").Value = 1312
$ws.Range("F5").Value = "Gammelgransskål"
$ws.Range("G5").Value = "Pseudographis pinicola"
$ws.Range("H5").Value = "(Nyl.) Rehm"
$ws.Range("R5").Value = 7069709

# Row 6
$ws.Range("A6").Value = 112388110
$ws.Range("B6").Value = 77650
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("L6").Value = ""
$ws.Range("Q6").Value = 557867
$ws.Range("R6").Value = 7069706

# Row 7
$ws.Range("A7").Value = 112388115
$ws.Range("B7").Value = 96735
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = ""
$ws.Range("Q7").Value = 557811
$ws.Range("R7").Value = 7069647

# Row 8
$ws.Range("A8").Value = 112388103
$ws.Range("B8").Value = 77650
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("L8").Value = ""
$ws.Range("Q8").Value = 557984
$ws.Range("R8").Value = 7069574

# Row 9
$ws.Range("A9").Value = 112388101
$ws.Range("B9").Value = 96735
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = ""
$ws.Range("Q9").Value = 557984
$ws.Range("R9").Value = 7069575
